$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 2).Value = -19.71713753387756
$ws.Cells.Item(2, 3).Value = 1.954702813604764
$ws.Cells.Item(2, 4).Value = -19.71713753387756
$ws.Cells.Item(2, 5).Value = -19.71713753387756
$ws.Cells.Item(2, 6).Value = -19.71713753387756
$ws.Cells.Item(2, 7).Value = -19.71713753387756
$ws.Cells.Item(2, 8).Value = -19.71713753387756
$ws.Cells.Item(2, 9).Value = -19.71713753387756
$ws.Cells.Item(2, 10).Value = -19.71713753387756
$ws.Cells.Item(2, 11).Value = -19.71713753387756
$ws.Cells.Item(3, 2).Value = -19.71713753387756
$ws.Cells.Item(3, 3).Value = -19.71713753387756
$ws.Cells.Item(3, 4).Value = -19.71713753387756
$ws.Cells.Item(3, 5).Value = -19.71713753387756
$ws.Cells.Item(3, 6).Value = -19.71713753387756
$ws.Cells.Item(3, 7).Value = -19.71713753387756
$ws.Cells.Item(3, 8).Value = -19.71713753387756
$ws.Cells.Item(3, 9).Value = 1.552108450027562
$ws.Cells.Item(3, 10).Value = -19.71713753387756
$ws.Cells.Item(3, 11).Value = -19.71713753387756
$ws.Cells.Item(4, 2).Value = -19.71713753387756
$ws.Cells.Item(4, 3).Value = 2.00178032945691
$ws.Cells.Item(4, 4).Value = 1.694634621884136
$ws.Cells.Item(4, 5).Value = -19.71713753387756
$ws.Cells.Item(4, 6).Value = 3.448362185677476
$ws.Cells.Item(4, 7).Value = -19.71713753387756
$ws.Cells.Item(4, 8).Value = 1.465793234619068
$ws.Cells.Item(4, 9).Value = -19.71713753387756
$ws.Cells.Item(4, 10).Value = -19.71713753387756
$ws.Cells.Item(4, 11).Value = -19.71713753387756
$ws.Cells.Item(5, 2).Value = -19.71713753387756
$ws.Cells.Item(5, 3).Value = 1.664895907855098
$ws.Cells.Item(5, 4).Value = -19.71713753387756
$ws.Cells.Item(5, 5).Value = -19.71713753387756
$ws.Cells.Item(5, 6).Value = -19.71713753387756
$ws.Cells.Item(5, 7).Value = 2.870144866393892
$ws.Cells.Item(5, 8).Value = -19.71713753387756
$ws.Cells.Item(5, 9).Value = -19.71713753387756
$ws.Cells.Item(5, 10).Value = -19.71713753387756
$ws.Cells.Item(5, 11).Value = -19.71713753387756
$ws.Cells.Item(6, 2).Value = -19.71713753387756
$ws.Cells.Item(6, 3).Value = -19.71713753387756
$ws.Cells.Item(6, 4).Value = -19.71713753387756
$ws.Cells.Item(6, 5).Value = -19.71713753387756
$ws.Cells.Item(6, 6).Value = -19.71713753387756
$ws.Cells.Item(6, 7).Value = -19.71713753387756
$ws.Cells.Item(6, 8).Value = -19.71713753387756
$ws.Cells.Item(6, 9).Value = -19.71713753387756
$ws.Cells.Item(6, 10).Value = -19.71713753387756
$ws.Cells.Item(6, 11).Value = -19.71713753387756
$ws.Cells.Item(7, 2).Value = 2.470083130333226
$ws.Cells.Item(7, 3).Value = -19.71713753387756
$ws.Cells.Item(7, 4).Value = -19.71713753387756
$ws.Cells.Item(7, 5).Value = -19.71713753387756
$ws.Cells.Item(7, 6).Value = -19.71713753387756
$ws.Cells.Item(7, 7).Value = -19.71713753387756
$ws.Cells.Item(7, 8).Value = -19.71713753387756
$ws.Cells.Item(7, 9).Value = -19.71713753387756
$ws.Cells.Item(7, 10).Value = -19.71713753387756
$ws.Cells.Item(7, 11).Value = -19.71713753387756
$ws.Cells.Item(8, 2).Value = -19.71713753387756
$ws.Cells.Item(8, 3).Value = -19.71713753387756
$ws.Cells.Item(8, 4).Value = -19.71713753387756
$ws.Cells.Item(8, 5).Value = 1.836312240633818
$ws.Cells.Item(8, 6).Value = -19.71713753387756
$ws.Cells.Item(8, 7).Value = -19.71713753387756
$ws.Cells.Item(8, 8).Value = -19.71713753387756
$ws.Cells.Item(8, 9).Value = -19.71713753387756
$ws.Cells.Item(8, 10).Value = -19.71713753387756
$ws.Cells.Item(8, 11).Value = -19.71713753387756
$ws.Cells.Item(9, 2).Value = 3.853918021437529
$ws.Cells.Item(9, 3).Value = -19.71713753387756
$ws.Cells.Item(9, 4).Value = -19.71713753387756
$ws.Cells.Item(9, 5).Value = -19.71713753387756
$ws.Cells.Item(9, 6).Value = -19.71713753387756
$ws.Cells.Item(9, 7).Value = -19.71713753387756
$ws.Cells.Item(9, 8).Value = -19.71713753387756
$ws.Cells.Item(9, 9).Value = -19.71713753387756
$ws.Cells.Item(9, 10).Value = -19.71713753387756
$ws.Cells.Item(9, 11).Value = -19.71713753387756
$ws.Cells.Item(10, 2).Value = -19.71713753387756
$ws.Cells.Item(10, 3).Value = -19.71713753387756
$ws.Cells.Item(10, 4).Value = -19.71713753387756
$ws.Cells.Item(10, 5).Value = -19.71713753387756
$ws.Cells.Item(10, 6).Value = -19.71713753387756
$ws.Cells.Item(10, 7).Value = -19.71713753387756
$ws.Cells.Item(10, 8).Value = -19.71713753387756
$ws.Cells.Item(10, 9).Value = 1.597552465065573
$ws.Cells.Item(10, 10).Value = -19.71713753387756
$ws.Cells.Item(10, 11).Value = 2.164190432481718
$ws.Cells.Item(11, 2).Value = -19.71713753387756
$ws.Cells.Item(11, 3).Value = -19.71713753387756
$ws.Cells.Item(11, 4).Value = -19.71713753387756
$ws.Cells.Item(11, 5).Value = 2.915243808483798
$ws.Cells.Item(11, 6).Value = -19.71713753387756
$ws.Cells.Item(11, 7).Value = 2.866114445660983
$ws.Cells.Item(11, 8).Value = -19.71713753387756
$ws.Cells.Item(11, 9).Value = -19.71713753387756
$ws.Cells.Item(11, 10).Value = -19.71713753387756
$ws.Cells.Item(11, 11).Value = 2.080690594510613
$ws.Cells.Item(12, 2).Value = -19.71713753387756
$ws.Cells.Item(12, 3).Value = -19.71713753387756
$ws.Cells.Item(12, 4).Value = -19.71713753387756
$ws.Cells.Item(12, 5).Value = -19.71713753387756
$ws.Cells.Item(12, 6).Value = -19.71713753387756
$ws.Cells.Item(12, 7).Value = -19.71713753387756
$ws.Cells.Item(12, 8).Value = -19.71713753387756
$ws.Cells.Item(12, 9).Value = -19.71713753387756
$ws.Cells.Item(12, 10).Value = -19.71713753387756
$ws.Cells.Item(12, 11).Value = -19.71713753387756
$ws.Cells.Item(13, 2).Value = -19.71713753387756
$ws.Cells.Item(13, 3).Value = -19.71713753387756
$ws.Cells.Item(13, 4).Value = -19.71713753387756
$ws.Cells.Item(13, 5).Value = 2.494698787909308
$ws.Cells.Item(13, 6).Value = -19.71713753387756
$ws.Cells.Item(13, 7).Value = -19.71713753387756
$ws.Cells.Item(13, 8).Value = -19.71713753387756
$ws.Cells.Item(13, 9).Value = -19.71713753387756
$ws.Cells.Item(13, 10).Value = -19.71713753387756
$ws.Cells.Item(13, 11).Value = 1.850631513804179
$ws.Cells.Item(14, 2).Value = -19.71713753387756
$ws.Cells.Item(14, 3).Value = -19.71713753387756
$ws.Cells.Item(14, 4).Value = 1.540600760144856
$ws.Cells.Item(14, 5).Value = -19.71713753387756
$ws.Cells.Item(14, 6).Value = -19.71713753387756
$ws.Cells.Item(14, 7).Value = -19.71713753387756
$ws.Cells.Item(14, 8).Value = -19.71713753387756
$ws.Cells.Item(14, 9).Value = -19.71713753387756
$ws.Cells.Item(14, 10).Value = -19.71713753387756
$ws.Cells.Item(14, 11).Value = 2.029285346453244
$ws.Cells.Item(15, 2).Value = -19.71713753387756
$ws.Cells.Item(15, 3).Value = -19.71713753387756
$ws.Cells.Item(15, 4).Value = 1.74434241301201
$ws.Cells.Item(15, 5).Value = -19.71713753387756
$ws.Cells.Item(15, 6).Value = -19.71713753387756
$ws.Cells.Item(15, 7).Value = -19.71713753387756
$ws.Cells.Item(15, 8).Value = -19.71713753387756
$ws.Cells.Item(15, 9).Value = -19.71713753387756
$ws.Cells.Item(15, 10).Value = -19.71713753387756
$ws.Cells.Item(15, 11).Value = -19.71713753387756
$ws.Cells.Item(16, 2).Value = -19.71713753387756
$ws.Cells.Item(16, 3).Value = -19.71713753387756
$ws.Cells.Item(16, 4).Value = -19.71713753387756
$ws.Cells.Item(16, 5).Value = -19.71713753387756
$ws.Cells.Item(16, 6).Value = -19.71713753387756
$ws.Cells.Item(16, 7).Value = -19.71713753387756
$ws.Cells.Item(16, 8).Value = -19.71713753387756
$ws.Cells.Item(16, 9).Value = -19.71713753387756
$ws.Cells.Item(16, 10).Value = -19.71713753387756
$ws.Cells.Item(16, 11).Value = -19.71713753387756
$ws.Cells.Item(17, 2).Value = -19.71713753387756
$ws.Cells.Item(17, 3).Value = 2.113491013619491
$ws.Cells.Item(17, 4).Value = 1.822521589710242
$ws.Cells.Item(17, 5).Value = -19.71713753387756
$ws.Cells.Item(17, 6).Value = -19.71713753387756
$ws.Cells.Item(17, 7).Value = -19.71713753387756
$ws.Cells.Item(17, 8).Value = 2.11908121435581
$ws.Cells.Item(17, 9).Value = 1.995942250226331
$ws.Cells.Item(17, 10).Value = 4.321926504699049
$ws.Cells.Item(17, 11).Value = -19.71713753387756
$ws.Cells.Item(18, 2).Value = -19.71713753387756
$ws.Cells.Item(18, 3).Value = -19.71713753387756
$ws.Cells.Item(18, 4).Value = -19.71713753387756
$ws.Cells.Item(18, 5).Value = -19.71713753387756
$ws.Cells.Item(18, 6).Value = -19.71713753387756
$ws.Cells.Item(18, 7).Value = -19.71713753387756
$ws.Cells.Item(18, 8).Value = 1.954486888550699
$ws.Cells.Item(18, 9).Value = 1.94152418027749
$ws.Cells.Item(18, 10).Value = -19.71713753387756
$ws.Cells.Item(18, 11).Value = -19.71713753387756
$ws.Cells.Item(19, 2).Value = -19.71713753387756
$ws.Cells.Item(19, 3).Value = -19.71713753387756
$ws.Cells.Item(19, 4).Value = 2.02635453269158
$ws.Cells.Item(19, 5).Value = -19.71713753387756
$ws.Cells.Item(19, 6).Value = -19.71713753387756
$ws.Cells.Item(19, 7).Value = -19.71713753387756
$ws.Cells.Item(19, 8).Value = 1.653767467050973
$ws.Cells.Item(19, 9).Value = 1.840730451753398
$ws.Cells.Item(19, 10).Value = -19.71713753387756
$ws.Cells.Item(19, 11).Value = -19.71713753387756
$ws.Cells.Item(20, 2).Value = -19.71713753387756
$ws.Cells.Item(20, 3).Value = 1.05851212752642
$ws.Cells.Item(20, 4).Value = 1.532430231959881
$ws.Cells.Item(20, 5).Value = -19.71713753387756
$ws.Cells.Item(20, 6).Value = 3.183336656971275
$ws.Cells.Item(20, 7).Value = -19.71713753387756
$ws.Cells.Item(20, 8).Value = 1.645808488813632
$ws.Cells.Item(20, 9).Value = 1.394757789885052
$ws.Cells.Item(20, 10).Value = -19.71713753387756
$ws.Cells.Item(20, 11).Value = 1.84764879019772
$ws.Cells.Item(21, 2).Value = -19.71713753387756
$ws.Cells.Item(21, 3).Value = 1.343477266094902
$ws.Cells.Item(21, 4).Value = -19.71713753387756
$ws.Cells.Item(21, 5).Value = 1.700157190783616
$ws.Cells.Item(21, 6).Value = -19.71713753387756
$ws.Cells.Item(21, 7).Value = 2.432319656726731
$ws.Cells.Item(21, 8).Value = 1.455283259977323
$ws.Cells.Item(21, 9).Value = -19.71713753387756
$ws.Cells.Item(21, 10).Value = -19.71713753387756
$ws.Cells.Item(21, 11).Value = -19.71713753387756
